$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: "object, it receives +10" -> "object, it receives +1"
#   Merge the "+1"/"0"/bookmark/" reward...reward" runs into a single run,
#   drop the (old) _GoBack bookmark, while keeping the "target " run before
#   and the " is received" run after untouched.  We bracket the edit zone
#   with temporary bookmarks so the engine's run-coalescing stays inside it.
# ---------------------------------------------------------------------------
$full = $d.Content.Text
$idxL = $full.IndexOf("object, it receives")
$idxR = $full.IndexOf(" is received")
$rL = $d.Range($idxL, $idxL)
$rR = $d.Range($idxR, $idxR)
$d.Bookmarks.Add("ZZ_TempL1", $rL) | Out-Null
$d.Bookmarks.Add("ZZ_TempR1", $rR) | Out-Null

$d.Bookmarks("_GoBack").Delete()

$d.Content.Find.Execute("receives +10", $true, $false, $false, $false, $false, $true, 1, $false, "receives +1", 2)

$d.Bookmarks("ZZ_TempL1").Delete()
$d.Bookmarks("ZZ_TempR1").Delete()

# ---------------------------------------------------------------------------
# Part 2: "receives a small reward" -> "receives a very small reward", with
#   the text split into three runs: "receives a " | "very " | "small reward
#   ...artifacts", and the _GoBack bookmark placed between "very " and
#   "small reward...".
# ---------------------------------------------------------------------------
$full2 = $d.Content.Text
$idxL2 = $full2.IndexOf("receives a small")
$idxR2 = $full2.IndexOf(" of the observation")
$rL2 = $d.Range($idxL2, $idxL2)
$rR2 = $d.Range($idxR2, $idxR2)
$d.Bookmarks.Add("ZZ_TempL2", $rL2) | Out-Null
$d.Bookmarks.Add("ZZ_TempR2", $rR2) | Out-Null

$d.Content.Find.Execute("receives a small reward", $true, $false, $false, $false, $false, $true, 1, $false, "receives a very small reward", 2)

$d.Bookmarks("ZZ_TempL2").Delete()
$d.Bookmarks("ZZ_TempR2").Delete()

# Split the merged run into "receives a " | "very " | _GoBack | "small reward..."
$full3 = $d.Content.Text
$idxSplit1 = $full3.IndexOf("very small reward")
$r1 = $d.Range($idxSplit1, $idxSplit1)
$d.Bookmarks.Add("ZZ_SplitPoint", $r1) | Out-Null

$full4 = $d.Content.Text
$idxSplit2 = $full4.IndexOf("small reward dependently")
$r2 = $d.Range($idxSplit2, $idxSplit2)
$d.Bookmarks.Add("_GoBack", $r2) | Out-Null

$d.Bookmarks("ZZ_SplitPoint").Delete()
